# Insert a new data row at row 201 (shifting the existing rows 201:235 down
# to 202:236) and populate it with the new "Poroto verde" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("201:201").Insert()

$ws.Range("A201").Value = 10
$ws.Range("B201").Value = "Vega Modelo de Temuco"
$ws.Range("C201").Value = "La Araucanía"
$ws.Range("D201").Value = 45218
$ws.Range("E201").Value = 9
$ws.Range("F201").Value = 100112031
$ws.Range("G201").Value = "Poroto verde"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 65
$ws.Range("K201").Value = 1800
$ws.Range("L201").Value = 1800
$ws.Range("M201").Value = 1800
$ws.Range("N201").Value = "$/saco 25 kilos"
$ws.Range("O201").Value = "Provincia de Limarí"
$ws.Range("P201").Value = 72
$ws.Range("Q201").Value = 25
$ws.Range("R201").Value = "Hortaliza"
